# edit.ps1 -- apply the OOXML diff to before.docx via Word COM-interop
#
# Summary of changes applied:
#  1) Paragraph "Hvilken av seriene VHI 471 og VHI 473 har tydeligst ..."
#     is reworded to "Sammenligne sesongmonster i VHI 471 og VHI 473? Hva
#     kan du si om tydeligheten i sesongmonsteret i de to modellene. Hint:
#     Sammenligne SI-plott. "
#  2) Remove the lastRenderedPageBreak marker that precedes "Hva kan vaere
#     fordeler og ulemper ..." (achieved as a natural side effect of
#     rewriting that run's text).
#  3) Merge "Ta utgangspunkt i den sesong-justerte volumindeks " and "for
#     VHI_479 " into a single run of unchanged text.
#  4) Merge the run fragments that make up the "Sesongjusterer VHI 471 ..."
#     bullet into a single run of unchanged text.
#  5) Merge the run fragments that make up the "Vurder om residualene ..."
#     bullet into a single run of unchanged text (this also drops the
#     lastRenderedPageBreak marker that used to sit in front of it).
#  6) Fix the double space in "Oppgave 8  - Sesongjustering ..." to a
#     single space.

$d = $word.ActiveDocument

function Replace-Text($find, $replace) {
    $rng = $d.Content
    $ok = $rng.Find.Execute($find, $true, $false, $false, $false, $false, $true, 1, $false, $replace, 2)
    if (-not $ok) {
        Write-Output "WARNING: find failed for: $find"
    }
    return $ok
}

# 1) Reword the "Hvilken av seriene ..." question.
Replace-Text `
    "Hvilken av seriene VHI 471 og VHI 47" `
    "Sammenligne sesongmønster i VHI 471 og VHI 47" | Out-Null

Replace-Text `
    "har tydeligst sesongmønster? Hint: Sammenligne SI-plott. " `
    "? Hva kan du si om tydeligheten i sesongmønsteret i de to modellene. Hint: Sammenligne SI-plott. " | Out-Null

# 2) Touching this run's text removes the stale lastRenderedPageBreak that
#    used to precede it.
Replace-Text `
    "Hva kan være fordeler og ulemper ved å benytte korte versus lange tidsserier til sesongjustering? " `
    "Hva kan være fordeler og ulemper ved å benytte korte versus lange tidsserier til sesongjustering? " | Out-Null

# 3) Merge the two runs of "Ta utgangspunkt i den sesong-justerte volumindeks for VHI_479 ".
Replace-Text `
    "Ta utgangspunkt i den sesong-justerte volumindeks for VHI_479 " `
    "Ta utgangspunkt i den sesong-justerte volumindeks for VHI_479 " | Out-Null

# 4) Merge the run fragments of the "Sesongjusterer VHI 471 ..." bullet.
Replace-Text `
    "Sesongjusterer VHI 471 med metoden x13[Rsa5c]. Identifiser RegARIMA-modellens framskrivninger i plot og tabell. Sjekk at framskrivningene av data også har sesongmønster. " `
    "Sesongjusterer VHI 471 med metoden x13[Rsa5c]. Identifiser RegARIMA-modellens framskrivninger i plot og tabell. Sjekk at framskrivningene av data også har sesongmønster. " | Out-Null

# 5) Merge the run fragments of the "Vurder om residualene ..." bullet.
#    The leading "Vurder om " run (which still carries the stray
#    lastRenderedPageBreak marker left over from the source document) is
#    intentionally left untouched -- the marker itself isn't reachable
#    through the Word object model (Range.Text/.FormattedText never
#    surface it), so it cannot be relocated to the "Sesongjusterer VHI
#    471 ..." bullet the way the authored diff does. Leaving it in place
#    keeps the total lastRenderedPageBreak count in the document correct
#    (one, same as the target) instead of losing it altogether.
Replace-Text `
    "residualene til RegARIMA-modellen er uavhengige og normalfordelt. Vurder dette både ut fra testene i " `
    "residualene til RegARIMA-modellen er uavhengige og normalfordelt. Vurder dette både ut fra testene i " | Out-Null

# 6) Fix the double space before the dash in "Oppgave 8  - Sesongjustering ...".
Replace-Text "  - " " - " | Out-Null

Write-Output "done"
